$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 11 (shifts the old rows 11-16 down to 12-17,
# extending the sheet's used range to row 17).
$ws.Rows("11:11").Insert()

# Populate the newly inserted row 11 with the "PLPremium" spec line.
$ws.Range("A11").Value = "PLPremium"
$ws.Range("G11").Value = 0.735
$ws.Range("H11").Value = 0.8

# Move/leave the active selection on H12, matching the end-of-task cursor
# position.
[void]$ws.Range("H12").Select()
